$d = $word.ActiveDocument
$rng = $d.Range(0, 10)
$d.Bookmarks.Add("TestBookmark", $rng)
Write-Output "done"
